$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Update cell C10 (row 10, col C) value from 18 to 100 as captured by the diff.
$ws.Range("C10").Value = 100
